$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.369850754737854
$ws.Range("B1").Value = 1.446391582489014
$ws.Range("C1").Value = 1.632881760597229
$ws.Range("D1").Value = 2.576910972595215
$ws.Range("E1").Value = 4.595821857452393
